$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder BOM rows -------------------------------------------------
# Current order (rows 2-12):
#   2 Feather OLED
#   3 0.1" breakaway header
#   4 Vertical RJ45
#   5 Reed Switch            <- to be removed
#   6 Microcontroller        <- move to row 2
#   7 Motor Driver           <- move to row 3
#   8 Magnet                 <- to be removed
#   9 47 uF 35V Capacitor
#  10 4 pin connector plug
#  11 4 pin connector socket
#  12 12" IDC cable
#
# Target order (rows 2-11):
#   2 Microcontroller
#   3 Motor Driver
#   4 Feather OLED
#   5 0.1" breakaway header
#   6 Vertical RJ45
#   7 47 uF 35V Capacitor
#   8 4 pin connector plug
#   9 4 pin connector socket
#  10 12" IDC cable
#  11 3.3v dc regulator      (new row, Description only)

# Move "Motor Driver" (row 7) up to row 2 (inserting copied cells shifts
# everything else down by one row, so "Microcontroller" ends up at row 7).
$ws.Rows("7").Copy()
$ws.Rows("2").Insert()

# Move "Microcontroller" (now at row 7) up to row 2 as well, pushing
# "Motor Driver" (currently at row 2) down to row 3.
$ws.Rows("7").Copy()
$ws.Rows("2").Insert()

# Re-apply the original cell formatting for the Microcontroller
# description cell (it used a distinct style from the rest of column B).
# Row 8 now holds the duplicate "Microcontroller" row created by the
# insert above (row 7 is the original, untouched "Reed Switch" row).
$ws.Range("B8").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# Remove the now-duplicated rows plus the discontinued parts:
#   7  Reed Switch (original row, never re-used)
#   8  Microcontroller (duplicate created by the insert above)
#   9  Motor Driver (duplicate created by the insert above)
#  10  Magnet (discontinued part)
$ws.Range("A7:A10").EntireRow.Delete()

# --- Add the new BOM line ---------------------------------------------
$ws.Range("B11").Value = "3.3v dc regulator"

# --- Restore the selection shown in the saved workbook -----------------
$ws.Range("B12").Select()
